$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K62").Value = 2995.8572
$ws.Range("M62").Value = -2371.8572
$ws.Range("H62").Value = 3007.7896
$ws.Range("I62").Value = 2995.8572
$ws.Range("J62").Value = 3014.75
$ws.Range("N62").Value = -4262.75
$ws.Range("L62").Value = 3014.75
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996
$ws.Range("H64").Value = 4220.5
$ws.Range("L65").Value = 15073.75
$ws.Range("J65").Value = 3014.75
$ws.Range("M65").Value = -11859.286
$ws.Range("K65").Value = 14979.286
$ws.Range("H65").Value = 3007.7896
$ws.Range("N65").Value = -21313.75
$ws.Range("I65").Value = 2995.8572
$ws.Range("L67").Value = 4500
$ws.Range("H67").Value = 4220.5
$ws.Range("J67").Value = 4500
$ws.Range("N67").Value = -6216
$ws.Range("M74").Value = -2788.875
$ws.Range("K74").Value = 3724.875
$ws.Range("I74").Value = 3724.875
$ws.Range("H74").Value = 58866.555
$ws.Range("H76").Value = 4316.25
$ws.Range("M76").Value = -4001.25
$ws.Range("I76").Value = 4316.25
$ws.Range("K76").Value = 4316.25
$ws.Range("I77").Value = 3724.875
$ws.Range("H77").Value = 58866.555
$ws.Range("M77").Value = -13944.375
$ws.Range("K77").Value = 18624.375
$ws.Range("H79").Value = 4316.25
$ws.Range("I79").Value = 4316.25
$ws.Range("M79").Value = -3224.25
$ws.Range("K79").Value = 4316.25
$ws.Range("H80").Value = 1343.375
$ws.Range("J80").Value = 2748.5
$ws.Range("K80").Value = 2625
$ws.Range("N80").Value = -10241.5
$ws.Range("L80").Value = 8245.5
$ws.Range("I80").Value = 875
$ws.Range("M80").Value = -1627
$ws.Range("K83").Value = 7875
$ws.Range("L83").Value = 24736.5
$ws.Range("I83").Value = 875
$ws.Range("M83").Value = -2883
$ws.Range("H83").Value = 1343.375
$ws.Range("J83").Value = 2748.5
$ws.Range("N83").Value = -34720.5
$ws.Range("L87").Value = 50000
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("J90").Value = 50000
$ws.Range("N90").Value = -162480
$ws.Range("L90").Value = 150000
$ws.Range("H90").Value = 50000
$ws.Range("M98").Value = 159.0588
$ws.Range("K98").Value = 1338.9412
$ws.Range("I98").Value = 1338.9412
$ws.Range("H98").Value = 1338.9412
$ws.Range("H112").Value = 2084972
$ws.Range("J112").Value = 1914.375
$ws.Range("I112").Value = 6251087.5
$ws.Range("K112").Value = 18753262.5
$ws.Range("L112").Value = 5743.125
$ws.Range("M112").Value = -18752154.5
$ws.Range("N112").Value = -7959.125
$ws.Range("M122").Value = -1566.8236
$ws.Range("I122").Value = 1338.9412
$ws.Range("H122").Value = 1338.9412
$ws.Range("K122").Value = 4016.8236
$ws.Range("H138").Value = 1958

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5636.852
$ws.Range("I32").Value = 5307.8
$ws.Range("K32").Value = 5307.8
$ws.Range("M32").Value = -5020.8
$ws.Range("N32").Value = -10324
$ws.Range("J32").Value = 9750
$ws.Range("L32").Value = 9750
$ws.Range("J63").Value = 3383.3333
$ws.Range("L63").Value = 3383.3333
$ws.Range("H63").Value = 1852.762
$ws.Range("N63").Value = -4755.3333
$ws.Range("K63").Value = 1240.5333
$ws.Range("I63").Value = 1240.5333
$ws.Range("M63").Value = -554.5333000000001
$ws.Range("H66").Value = 1852.762
$ws.Range("J66").Value = 3383.3333
$ws.Range("I66").Value = 1240.5333
$ws.Range("K66").Value = 6202.6665
$ws.Range("L66").Value = 16916.6665
$ws.Range("M66").Value = -2770.6665
$ws.Range("N66").Value = -23780.6665
$ws.Range("K97").Value = 1087.4
$ws.Range("H97").Value = 1087.4
$ws.Range("M97").Value = -591.4000000000001
$ws.Range("I97").Value = 1087.4
$ws.Range("L130").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("N130").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K26").Value = 17500
$ws.Range("H26").Value = 17500
$ws.Range("I26").Value = 17500
$ws.Range("M26").Value = -17208
$ws.Range("K86").Value = 1896
$ws.Range("I86").Value = 1896
$ws.Range("H86").Value = 2511.5334
$ws.Range("M86").Value = -773
$ws.Range("H89").Value = 2511.5334
$ws.Range("K89").Value = 9480
$ws.Range("M89").Value = -3864
$ws.Range("I89").Value = 1896
$ws.Range("I99").Value = 2200
$ws.Range("M99").Value = -702
$ws.Range("K99").Value = 2200
$ws.Range("H99").Value = 2200
$ws.Range("K134").Value = 25470.669
$ws.Range("I134").Value = 8490.223
$ws.Range("L134").Value = 48000
$ws.Range("M134").Value = -22935.669
$ws.Range("H134").Value = 9241.200000000001
$ws.Range("J134").Value = 16000
$ws.Range("N134").Value = -53070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 1000
$ws.Range("K134").Value = 11407.5
$ws.Range("I134").Value = 3802.5
$ws.Range("M134").Value = -8872.5
$ws.Range("H134").Value = 3883.2666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4890.4443
$ws.Range("I81").Value = 4749.75
$ws.Range("N81").Value = -17255
$ws.Range("K81").Value = 14249.25
$ws.Range("L81").Value = 15009
$ws.Range("M81").Value = -13126.25
$ws.Range("J81").Value = 5003
$ws.Range("I84").Value = 4749.75
$ws.Range("H84").Value = 4890.4443
$ws.Range("J84").Value = 5003
$ws.Range("N84").Value = -56259
$ws.Range("M84").Value = -37131.75
$ws.Range("L84").Value = 45027
$ws.Range("K84").Value = 42747.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4486.4
$ws.Range("J80").Value = 6100
$ws.Range("K80").Value = 4083
$ws.Range("N80").Value = -8096
$ws.Range("L80").Value = 6100
$ws.Range("I80").Value = 4083
$ws.Range("M80").Value = -3085
$ws.Range("K83").Value = 20415
$ws.Range("L83").Value = 30500
$ws.Range("I83").Value = 4083
$ws.Range("M83").Value = -15423
$ws.Range("H83").Value = 4486.4
$ws.Range("J83").Value = 6100
$ws.Range("N83").Value = -40484
$ws.Range("K102").Value = 1197.2858
$ws.Range("H102").Value = 1197.2858
$ws.Range("I102").Value = 1197.2858
$ws.Range("M102").Value = 424.7141999999999
$ws.Range("K126").Value = 5908.9998
$ws.Range("H126").Value = 1969.6666
$ws.Range("M126").Value = -3438.9998
$ws.Range("I126").Value = 1969.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3718.318
$ws.Range("I46").Value = 3545.3635
$ws.Range("K46").Value = 3545.3635
$ws.Range("M46").Value = -3357.3635
$ws.Range("H55").Value = 399.75
$ws.Range("I55").Value = 200
$ws.Range("M55").Value = -27
$ws.Range("K55").Value = 200
$ws.Range("J55").Value = 466.33334
$ws.Range("L55").Value = 466.33334
$ws.Range("N55").Value = -812.33334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 20999.75
$ws.Range("J80").Value = 20999.75
$ws.Range("N80").Value = -22995.75
$ws.Range("L80").Value = 20999.75
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("L83").Value = 62999.25
$ws.Range("H83").Value = 20999.75
$ws.Range("J83").Value = 20999.75
$ws.Range("N83").Value = -72983.25
$ws.Range("I84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("K84").Value = 0
$ws.Range("H107").Value = 900
$ws.Range("K107").Value = 2700
$ws.Range("M107").Value = -780
$ws.Range("I107").Value = 900
$ws.Range("M122").Value = -586.4287000000004
$ws.Range("I122").Value = 1012.1429
$ws.Range("J122").Value = 905
$ws.Range("H122").Value = 998.75
$ws.Range("N122").Value = -7615
$ws.Range("K122").Value = 3036.4287
$ws.Range("L122").Value = 2715
